# Fixed not displaying correct old efficiency
# 1) Update file path strings: change "/" separators to "\" and
#    "surveys" folder to "Output\Output Files"
# 2) Update several "Old Efficiency" (and matching "New Efficiency")
#    values that were incorrect.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Fix file path text in column A (old "surveys" folder -> new "Output\Output Files" folder) ---
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*surveys\*") {
        $newVal = $val -replace "/", "\"
        $newVal = $newVal -replace "\\surveys\\", "\Output\Output Files\"
        $cell.Value2 = $newVal
    }
}

# --- 2) Fix incorrect Old Efficiency (and matching New Efficiency) values ---
$ws.Cells.Item(14, 5).Value = 0.13

$ws.Cells.Item(18, 5).Value = 0.13
$ws.Cells.Item(19, 5).Value = 0.13
$ws.Cells.Item(20, 5).Value = 0.13
$ws.Cells.Item(21, 5).Value = 0.13
$ws.Cells.Item(22, 5).Value = 0.13

$ws.Cells.Item(49, 5).Value = 0.01
$ws.Cells.Item(49, 6).Value = 0.01
$ws.Cells.Item(50, 5).Value = 0.01
$ws.Cells.Item(50, 6).Value = 0.01
$ws.Cells.Item(65, 5).Value = 0.01
$ws.Cells.Item(65, 6).Value = 0.01
